$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 22 (Empty case)
$ws.Range("D22").Value = 102
$ws.Range("I22").Value = 38

# Row 23 (NumberOnly case)
$ws.Range("D23").Value = 91
$ws.Range("I23").Value = 54

# Row 24
$ws.Range("D24").Value = 143
$ws.Range("I24").Value = 59

# Row 25
$ws.Range("D25").Value = 480
$ws.Range("E25").Value = 1133048
$ws.Range("I25").Value = 107

# Row 26
$ws.Range("D26").Value = 597
$ws.Range("I26").Value = 118

# Row 27
$ws.Range("D27").Value = 2552
$ws.Range("E27").Value = 6954199
$ws.Range("I27").Value = 351

# Row 28
$ws.Range("D28").Value = 352
$ws.Range("E28").Value = 609594
$ws.Range("I28").Value = 159

# Row 29
$ws.Range("D29").Value = 1342
$ws.Range("E29").Value = 3305395
$ws.Range("I29").Value = 370

# Row 30
$ws.Range("D30").Value = 323
$ws.Range("E30").Value = 804816

# Row 31
$ws.Range("D31").Value = 695
$ws.Range("E31").Value = 2094014
$ws.Range("I31").Value = 284

# Row 32
$ws.Range("D32").Value = 8567
$ws.Range("E32").Value = 9360871
$ws.Range("I32").Value = 622

# Row 33
$ws.Range("D33").Value = 717
$ws.Range("E33").Value = 1515846
$ws.Range("I33").Value = 250

# Clear the lingering selection left over from the previous edit session
# (source file had D13:U13 selected; reset to the sheet's home cell).
[void]$ws.Range("A1").Select()

$wb.Save()
